$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 165; this shifts old rows 165-192 down to 166-193.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new data point.
$ws.Cells.Item(165, 1).Value = 5
$ws.Cells.Item(165, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(165, 3).Value = "Maule"
$ws.Cells.Item(165, 4).Value = 44491
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(165, 5).Value = 7
$ws.Cells.Item(165, 6).Value = 100114014
$ws.Cells.Item(165, 7).Value = "Betarraga"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 4000
$ws.Cells.Item(165, 11).Value = 650
$ws.Cells.Item(165, 12).Value = 650
$ws.Cells.Item(165, 13).Value = 650
$ws.Cells.Item(165, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 130
$ws.Cells.Item(165, 17).Value = 5
$ws.Cells.Item(165, 18).Value = "Hortaliza"
